$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the vessel type / description / weight / counts for row 3
# (261k02), keeping only the benefit code.
$ws.Range("B3:G3").ClearContents()

# Add a new row (row 4) for benefit code 262k01.
$ws.Range("A4").Value = "262k01"
$ws.Range("B4").Value = "steel"
$ws.Range("C4").Value = "bucket"
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

$ws.Range("H4").Select()
